# Re-procesar los datos con las nuevas dimensiones curadas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columna C (municipio-nombre): ya no es iaest-measure, ahora es sdmx-dimension:refArea / dim / URI-Municipio
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("C3").Value = "dim"
$ws.Range("C4").Value = "URI-Municipio"

# Columna H (porcentaje-de-sau-en-propiedad-del-titular): ya no es dimension, ahora es iaest-measure / medida / xsd:int
$ws.Range("H2").Value = "iaest-measure:porcentaje-de-sau-en-propiedad-del-titular"
$ws.Range("H3").Value = "medida"
$ws.Range("H4").Value = "xsd:int"

# Ya no requiere fichero de mapeo (ya no es una dimension curada por mapeo)
$ws.Range("H5").Clear()
